$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 4 ("register.tlp" row) to make room
# for the new es-es / fr-fr / ru-ru account rows for register.php.
$ws.Rows.Item(4).Resize(3).Insert()

# register.php additional language rows (rows 4-6)
$ws.Range("A4").Value = "register.php"
$ws.Range("B4").Value = "upload\catalog\language\es-es\account"
$ws.Range("C4").Value = "Jesus Martinez"

$ws.Range("A5").Value = "register.php"
$ws.Range("B5").Value = "upload\catalog\language\fr-fr\account"
$ws.Range("C5").Value = "Jesus Martinez"

$ws.Range("A6").Value = "register.php"
$ws.Range("B6").Value = "upload\catalog\language\ru-ru\account"
$ws.Range("C6").Value = "Jesus Martinez"

# register.tlp row (now at row 7, was row 4 previously)
$ws.Range("A7").Value = "register.tlp"
$ws.Range("B7").Value = "upload\catalog\view\theme\ArtsBoutiqueLvovna\template\account"
$ws.Range("C7").Value = "Jesus Martinez"

# customer.php rows (new rows 8-9)
$ws.Range("B8").Value = "upload\catalog\model\account\"
$ws.Range("A8").Value = "customer.php"
$ws.Range("C8").Value = "Jesus Martinez"

$ws.Range("A9").Value = "customer.php"
$ws.Range("B9").Value = "upload\admin\model\customer\"
$ws.Range("C9").Value = "Jesus Martinez"

$ws.Range("C10").Select()
